$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet by copying the "2022-Q2" sheet
#    (so it inherits identical formatting/layout), inserting it right
#    before "2022-Q2" (i.e. right after "总计").
# ------------------------------------------------------------------
$sheetQ2 = $wb.Worksheets.Item("2022-Q2")
$sheetQ2Index = $sheetQ2.Index
$sheetQ2.Copy($sheetQ2)
$newSheet = $wb.Worksheets.Item($sheetQ2Index)
$newSheet.Name = "2022-Q3"

# Extend the table from 5 data rows (2..6) to 7 data rows (2..8),
# copying the formatting of the last existing data row down.
$newSheet.Range("A6:H6").Copy()
$newSheet.Range("A7:H8").PasteSpecial(-4122)

# Fill in the new fund-holding data for 2022-Q3.
$rows = @(
    @(0, "008269", "大成睿享混合A", "19.80", "66.91", "6.04", "1.1959", 1),
    @(1, "090013", "大成竞争优势混合", "6.88", "61.00", "6.27", "0.4314", 1),
    @(2, "013463", "大成致远优势一年持有期混合A", "3.65", "60.88", "8.71", "0.3179", 1),
    @(3, "008270", "大成睿享混合C", "4.02", "66.91", "6.04", "0.2428", 1),
    @(4, "014094", "南方誉盈一年持有混合A", "11.86", "24.33", "0.69", "0.0818", 10),
    @(5, "013464", "大成致远优势一年持有期混合C", "0.17", "60.88", "8.71", "0.0148", 1),
    @(6, "014095", "南方誉盈一年持有混合C", "0.84", "24.33", "0.69", "0.0058", 10)
)

$newSheet.Range("B2:G8").NumberFormat = "@"

$r = 2
foreach ($row in $rows) {
    $newSheet.Range("A$r").Value = $row[0]
    $newSheet.Range("B$r").Value = $row[1]
    $newSheet.Range("C$r").Value = $row[2]
    $newSheet.Range("D$r").Value = $row[3]
    $newSheet.Range("E$r").Value = $row[4]
    $newSheet.Range("F$r").Value = $row[5]
    $newSheet.Range("G$r").Value = $row[6]
    $newSheet.Range("H$r").Value = $row[7]
    $r = $r + 1
}

$newSheet.Range("B2:G8").Style = "Normal"

# ------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a new row for 2022-Q3
#    above the 2022-Q2 row, shifting the remaining rows down and
#    renumbering the index column.
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("总计")
$ws.Rows.Item(2).Insert()
$ws.Range("B2:D2").Style = "Normal"

$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "2022-Q3"
$ws.Range("C2").Value = 7
$ws.Range("D2").Value = 2.29

$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
